$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rat Colony")

# --- Step 1: fill in missing Weight_g (column C) values for existing rows 668-685 ---
$ws.Cells.Item(668, 3).Value = 399
$ws.Cells.Item(669, 3).Value = 398
$ws.Cells.Item(670, 3).Value = 374
$ws.Cells.Item(671, 3).Value = 407
$ws.Cells.Item(672, 3).Value = 421
$ws.Cells.Item(673, 3).Value = 424
$ws.Cells.Item(674, 3).Value = 376
$ws.Cells.Item(675, 3).Value = 421
$ws.Cells.Item(676, 3).Value = 429
$ws.Cells.Item(677, 3).Value = 407
$ws.Cells.Item(678, 3).Value = 409
$ws.Cells.Item(679, 3).Value = 377
$ws.Cells.Item(680, 3).Value = 422
$ws.Cells.Item(681, 3).Value = 431
$ws.Cells.Item(682, 3).Value = 429
$ws.Cells.Item(683, 3).Value = 382
$ws.Cells.Item(684, 3).Value = 425
$ws.Cells.Item(685, 3).Value = 433

# --- Step 2: append new rows 686-766 (9 new days of rat colony observations) ---
$ws.Cells.Item(686, 1).Value = 1
$ws.Cells.Item(686, 2).Formula = "=B677+1"
$ws.Cells.Item(686, 3).Value = 409
$ws.Cells.Item(686, 4).Value = 1
$ws.Cells.Item(686, 5).Value = "Black tipped tail"
$ws.Cells.Item(686, 6).Value = "BTT"
$ws.Cells.Item(687, 1).Value = 2
$ws.Cells.Item(687, 2).Formula = "=B686"
$ws.Cells.Item(687, 3).Value = 410
$ws.Cells.Item(687, 4).Value = 1
$ws.Cells.Item(687, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(687, 6).Value = "BDLRL"
$ws.Cells.Item(688, 1).Value = 3
$ws.Cells.Item(688, 2).Formula = "=B679+1"
$ws.Cells.Item(688, 3).Value = 379
$ws.Cells.Item(688, 4).Value = 1
$ws.Cells.Item(688, 5).Value = "White rear legs"
$ws.Cells.Item(688, 6).Value = "WRL"
$ws.Cells.Item(689, 1).Value = 4
$ws.Cells.Item(689, 2).Formula = "=B688"
$ws.Cells.Item(689, 3).Value = 410
$ws.Cells.Item(689, 4).Value = 2
$ws.Cells.Item(689, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(689, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(690, 1).Value = 5
$ws.Cells.Item(690, 2).Formula = "=B681+1"
$ws.Cells.Item(690, 3).Value = 427
$ws.Cells.Item(690, 4).Value = 2
$ws.Cells.Item(690, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(690, 6).Value = "1/2WT/WU"
$ws.Cells.Item(691, 1).Value = 6
$ws.Cells.Item(691, 2).Formula = "=B682+1"
$ws.Cells.Item(691, 3).Value = 428
$ws.Cells.Item(691, 4).Value = 2
$ws.Cells.Item(691, 5).Value = "Small white tip tail"
$ws.Cells.Item(691, 6).Value = "SWTT"
$ws.Cells.Item(692, 1).Value = 7
$ws.Cells.Item(692, 2).Formula = "=B691"
$ws.Cells.Item(692, 3).Value = 383
$ws.Cells.Item(692, 4).Value = 3
$ws.Cells.Item(692, 5).Value = "White tail"
$ws.Cells.Item(692, 6).Value = "WT"
$ws.Cells.Item(693, 1).Value = 8
$ws.Cells.Item(693, 2).Formula = "=B684+1"
$ws.Cells.Item(693, 3).Value = 426
$ws.Cells.Item(693, 4).Value = 3
$ws.Cells.Item(693, 5).Value = "Half white tail"
$ws.Cells.Item(693, 6).Value = "1/2WT"
$ws.Cells.Item(694, 1).Value = 9
$ws.Cells.Item(694, 2).Formula = "=B693"
$ws.Cells.Item(694, 3).Value = 438
$ws.Cells.Item(694, 4).Value = 3
$ws.Cells.Item(694, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(694, 6).Value = "1/4WTT"
$ws.Cells.Item(695, 1).Value = 1
$ws.Cells.Item(695, 2).Formula = "=B686+1"
$ws.Cells.Item(695, 3).Value = 418
$ws.Cells.Item(695, 4).Value = 1
$ws.Cells.Item(695, 5).Value = "Black tipped tail"
$ws.Cells.Item(695, 6).Value = "BTT"
$ws.Cells.Item(696, 1).Value = 2
$ws.Cells.Item(696, 2).Formula = "=B695"
$ws.Cells.Item(696, 3).Value = 423
$ws.Cells.Item(696, 4).Value = 1
$ws.Cells.Item(696, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(696, 6).Value = "BDLRL"
$ws.Cells.Item(697, 1).Value = 3
$ws.Cells.Item(697, 2).Formula = "=B688+1"
$ws.Cells.Item(697, 3).Value = 393
$ws.Cells.Item(697, 4).Value = 1
$ws.Cells.Item(697, 5).Value = "White rear legs"
$ws.Cells.Item(697, 6).Value = "WRL"
$ws.Cells.Item(698, 1).Value = 4
$ws.Cells.Item(698, 2).Formula = "=B697"
$ws.Cells.Item(698, 3).Value = 429
$ws.Cells.Item(698, 4).Value = 2
$ws.Cells.Item(698, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(698, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(699, 1).Value = 5
$ws.Cells.Item(699, 2).Formula = "=B690+1"
$ws.Cells.Item(699, 3).Value = 442
$ws.Cells.Item(699, 4).Value = 2
$ws.Cells.Item(699, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(699, 6).Value = "1/2WT/WU"
$ws.Cells.Item(700, 1).Value = 6
$ws.Cells.Item(700, 2).Formula = "=B691+1"
$ws.Cells.Item(700, 3).Value = 436
$ws.Cells.Item(700, 4).Value = 2
$ws.Cells.Item(700, 5).Value = "Small white tip tail"
$ws.Cells.Item(700, 6).Value = "SWTT"
$ws.Cells.Item(701, 1).Value = 7
$ws.Cells.Item(701, 2).Formula = "=B700"
$ws.Cells.Item(701, 3).Value = 391
$ws.Cells.Item(701, 4).Value = 3
$ws.Cells.Item(701, 5).Value = "White tail"
$ws.Cells.Item(701, 6).Value = "WT"
$ws.Cells.Item(702, 1).Value = 8
$ws.Cells.Item(702, 2).Formula = "=B693+1"
$ws.Cells.Item(702, 3).Value = 434
$ws.Cells.Item(702, 4).Value = 3
$ws.Cells.Item(702, 5).Value = "Half white tail"
$ws.Cells.Item(702, 6).Value = "1/2WT"
$ws.Cells.Item(703, 1).Value = 9
$ws.Cells.Item(703, 2).Formula = "=B702"
$ws.Cells.Item(703, 3).Value = 443
$ws.Cells.Item(703, 4).Value = 3
$ws.Cells.Item(703, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(703, 6).Value = "1/4WTT"
$ws.Cells.Item(704, 1).Value = 1
$ws.Cells.Item(704, 2).Formula = "=B695+1"
$ws.Cells.Item(704, 3).Value = 411
$ws.Cells.Item(704, 4).Value = 1
$ws.Cells.Item(704, 5).Value = "Black tipped tail"
$ws.Cells.Item(704, 6).Value = "BTT"
$ws.Cells.Item(705, 1).Value = 2
$ws.Cells.Item(705, 2).Formula = "=B704"
$ws.Cells.Item(705, 3).Value = 421
$ws.Cells.Item(705, 4).Value = 1
$ws.Cells.Item(705, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(705, 6).Value = "BDLRL"
$ws.Cells.Item(706, 1).Value = 3
$ws.Cells.Item(706, 2).Formula = "=B697+1"
$ws.Cells.Item(706, 3).Value = 385
$ws.Cells.Item(706, 4).Value = 1
$ws.Cells.Item(706, 5).Value = "White rear legs"
$ws.Cells.Item(706, 6).Value = "WRL"
$ws.Cells.Item(707, 1).Value = 4
$ws.Cells.Item(707, 2).Formula = "=B706"
$ws.Cells.Item(707, 3).Value = 423
$ws.Cells.Item(707, 4).Value = 2
$ws.Cells.Item(707, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(707, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(708, 1).Value = 5
$ws.Cells.Item(708, 2).Formula = "=B699+1"
$ws.Cells.Item(708, 3).Value = 434
$ws.Cells.Item(708, 4).Value = 2
$ws.Cells.Item(708, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(708, 6).Value = "1/2WT/WU"
$ws.Cells.Item(709, 1).Value = 6
$ws.Cells.Item(709, 2).Formula = "=B700+1"
$ws.Cells.Item(709, 3).Value = 441
$ws.Cells.Item(709, 4).Value = 2
$ws.Cells.Item(709, 5).Value = "Small white tip tail"
$ws.Cells.Item(709, 6).Value = "SWTT"
$ws.Cells.Item(710, 1).Value = 7
$ws.Cells.Item(710, 2).Formula = "=B709"
$ws.Cells.Item(710, 3).Value = 385
$ws.Cells.Item(710, 4).Value = 3
$ws.Cells.Item(710, 5).Value = "White tail"
$ws.Cells.Item(710, 6).Value = "WT"
$ws.Cells.Item(711, 1).Value = 8
$ws.Cells.Item(711, 2).Formula = "=B702+1"
$ws.Cells.Item(711, 3).Value = 433
$ws.Cells.Item(711, 4).Value = 3
$ws.Cells.Item(711, 5).Value = "Half white tail"
$ws.Cells.Item(711, 6).Value = "1/2WT"
$ws.Cells.Item(712, 1).Value = 9
$ws.Cells.Item(712, 2).Formula = "=B711"
$ws.Cells.Item(712, 3).Value = 439
$ws.Cells.Item(712, 4).Value = 3
$ws.Cells.Item(712, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(712, 6).Value = "1/4WTT"
$ws.Cells.Item(713, 1).Value = 1
$ws.Cells.Item(713, 2).Formula = "=B704+1"
$ws.Cells.Item(713, 4).Value = 1
$ws.Cells.Item(713, 5).Value = "Black tipped tail"
$ws.Cells.Item(713, 6).Value = "BTT"
$ws.Cells.Item(714, 1).Value = 2
$ws.Cells.Item(714, 2).Formula = "=B713"
$ws.Cells.Item(714, 4).Value = 1
$ws.Cells.Item(714, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(714, 6).Value = "BDLRL"
$ws.Cells.Item(715, 1).Value = 3
$ws.Cells.Item(715, 2).Formula = "=B706+1"
$ws.Cells.Item(715, 4).Value = 1
$ws.Cells.Item(715, 5).Value = "White rear legs"
$ws.Cells.Item(715, 6).Value = "WRL"
$ws.Cells.Item(716, 1).Value = 4
$ws.Cells.Item(716, 2).Formula = "=B715"
$ws.Cells.Item(716, 4).Value = 2
$ws.Cells.Item(716, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(716, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(717, 1).Value = 5
$ws.Cells.Item(717, 2).Formula = "=B708+1"
$ws.Cells.Item(717, 4).Value = 2
$ws.Cells.Item(717, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(717, 6).Value = "1/2WT/WU"
$ws.Cells.Item(718, 1).Value = 6
$ws.Cells.Item(718, 2).Formula = "=B709+1"
$ws.Cells.Item(718, 4).Value = 2
$ws.Cells.Item(718, 5).Value = "Small white tip tail"
$ws.Cells.Item(718, 6).Value = "SWTT"
$ws.Cells.Item(719, 1).Value = 7
$ws.Cells.Item(719, 2).Formula = "=B718"
$ws.Cells.Item(719, 4).Value = 3
$ws.Cells.Item(719, 5).Value = "White tail"
$ws.Cells.Item(719, 6).Value = "WT"
$ws.Cells.Item(720, 1).Value = 8
$ws.Cells.Item(720, 2).Formula = "=B711+1"
$ws.Cells.Item(720, 4).Value = 3
$ws.Cells.Item(720, 5).Value = "Half white tail"
$ws.Cells.Item(720, 6).Value = "1/2WT"
$ws.Cells.Item(721, 1).Value = 9
$ws.Cells.Item(721, 2).Formula = "=B720"
$ws.Cells.Item(721, 4).Value = 3
$ws.Cells.Item(721, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(721, 6).Value = "1/4WTT"
$ws.Cells.Item(722, 1).Value = 1
$ws.Cells.Item(722, 2).Formula = "=B713+1"
$ws.Cells.Item(722, 4).Value = 1
$ws.Cells.Item(722, 5).Value = "Black tipped tail"
$ws.Cells.Item(722, 6).Value = "BTT"
$ws.Cells.Item(723, 1).Value = 2
$ws.Cells.Item(723, 2).Formula = "=B722"
$ws.Cells.Item(723, 4).Value = 1
$ws.Cells.Item(723, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(723, 6).Value = "BDLRL"
$ws.Cells.Item(724, 1).Value = 3
$ws.Cells.Item(724, 2).Formula = "=B715+1"
$ws.Cells.Item(724, 4).Value = 1
$ws.Cells.Item(724, 5).Value = "White rear legs"
$ws.Cells.Item(724, 6).Value = "WRL"
$ws.Cells.Item(725, 1).Value = 4
$ws.Cells.Item(725, 2).Formula = "=B724"
$ws.Cells.Item(725, 4).Value = 2
$ws.Cells.Item(725, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(725, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(726, 1).Value = 5
$ws.Cells.Item(726, 2).Formula = "=B717+1"
$ws.Cells.Item(726, 4).Value = 2
$ws.Cells.Item(726, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(726, 6).Value = "1/2WT/WU"
$ws.Cells.Item(727, 1).Value = 6
$ws.Cells.Item(727, 2).Formula = "=B718+1"
$ws.Cells.Item(727, 4).Value = 2
$ws.Cells.Item(727, 5).Value = "Small white tip tail"
$ws.Cells.Item(727, 6).Value = "SWTT"
$ws.Cells.Item(728, 1).Value = 7
$ws.Cells.Item(728, 2).Formula = "=B727"
$ws.Cells.Item(728, 4).Value = 3
$ws.Cells.Item(728, 5).Value = "White tail"
$ws.Cells.Item(728, 6).Value = "WT"
$ws.Cells.Item(729, 1).Value = 8
$ws.Cells.Item(729, 2).Formula = "=B720+1"
$ws.Cells.Item(729, 4).Value = 3
$ws.Cells.Item(729, 5).Value = "Half white tail"
$ws.Cells.Item(729, 6).Value = "1/2WT"
$ws.Cells.Item(730, 1).Value = 9
$ws.Cells.Item(730, 2).Formula = "=B729"
$ws.Cells.Item(730, 4).Value = 3
$ws.Cells.Item(730, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(730, 6).Value = "1/4WTT"
$ws.Cells.Item(731, 1).Value = 1
$ws.Cells.Item(731, 2).Formula = "=B722+1"
$ws.Cells.Item(731, 4).Value = 1
$ws.Cells.Item(731, 5).Value = "Black tipped tail"
$ws.Cells.Item(731, 6).Value = "BTT"
$ws.Cells.Item(732, 1).Value = 2
$ws.Cells.Item(732, 2).Formula = "=B731"
$ws.Cells.Item(732, 4).Value = 1
$ws.Cells.Item(732, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(732, 6).Value = "BDLRL"
$ws.Cells.Item(733, 1).Value = 3
$ws.Cells.Item(733, 2).Formula = "=B724+1"
$ws.Cells.Item(733, 4).Value = 1
$ws.Cells.Item(733, 5).Value = "White rear legs"
$ws.Cells.Item(733, 6).Value = "WRL"
$ws.Cells.Item(734, 1).Value = 4
$ws.Cells.Item(734, 2).Formula = "=B733"
$ws.Cells.Item(734, 4).Value = 2
$ws.Cells.Item(734, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(734, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(735, 1).Value = 5
$ws.Cells.Item(735, 2).Formula = "=B726+1"
$ws.Cells.Item(735, 4).Value = 2
$ws.Cells.Item(735, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(735, 6).Value = "1/2WT/WU"
$ws.Cells.Item(736, 1).Value = 6
$ws.Cells.Item(736, 2).Formula = "=B727+1"
$ws.Cells.Item(736, 4).Value = 2
$ws.Cells.Item(736, 5).Value = "Small white tip tail"
$ws.Cells.Item(736, 6).Value = "SWTT"
$ws.Cells.Item(737, 1).Value = 7
$ws.Cells.Item(737, 2).Formula = "=B736"
$ws.Cells.Item(737, 4).Value = 3
$ws.Cells.Item(737, 5).Value = "White tail"
$ws.Cells.Item(737, 6).Value = "WT"
$ws.Cells.Item(738, 1).Value = 8
$ws.Cells.Item(738, 2).Formula = "=B729+1"
$ws.Cells.Item(738, 4).Value = 3
$ws.Cells.Item(738, 5).Value = "Half white tail"
$ws.Cells.Item(738, 6).Value = "1/2WT"
$ws.Cells.Item(739, 1).Value = 9
$ws.Cells.Item(739, 2).Formula = "=B738"
$ws.Cells.Item(739, 4).Value = 3
$ws.Cells.Item(739, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(739, 6).Value = "1/4WTT"
$ws.Cells.Item(740, 1).Value = 1
$ws.Cells.Item(740, 2).Formula = "=B731+1"
$ws.Cells.Item(740, 4).Value = 1
$ws.Cells.Item(740, 5).Value = "Black tipped tail"
$ws.Cells.Item(740, 6).Value = "BTT"
$ws.Cells.Item(741, 1).Value = 2
$ws.Cells.Item(741, 2).Formula = "=B740"
$ws.Cells.Item(741, 4).Value = 1
$ws.Cells.Item(741, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(741, 6).Value = "BDLRL"
$ws.Cells.Item(742, 1).Value = 3
$ws.Cells.Item(742, 2).Formula = "=B733+1"
$ws.Cells.Item(742, 4).Value = 1
$ws.Cells.Item(742, 5).Value = "White rear legs"
$ws.Cells.Item(742, 6).Value = "WRL"
$ws.Cells.Item(743, 1).Value = 4
$ws.Cells.Item(743, 2).Formula = "=B742"
$ws.Cells.Item(743, 4).Value = 2
$ws.Cells.Item(743, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(743, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(744, 1).Value = 5
$ws.Cells.Item(744, 2).Formula = "=B735+1"
$ws.Cells.Item(744, 4).Value = 2
$ws.Cells.Item(744, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(744, 6).Value = "1/2WT/WU"
$ws.Cells.Item(745, 1).Value = 6
$ws.Cells.Item(745, 2).Formula = "=B736+1"
$ws.Cells.Item(745, 4).Value = 2
$ws.Cells.Item(745, 5).Value = "Small white tip tail"
$ws.Cells.Item(745, 6).Value = "SWTT"
$ws.Cells.Item(746, 1).Value = 7
$ws.Cells.Item(746, 2).Formula = "=B745"
$ws.Cells.Item(746, 4).Value = 3
$ws.Cells.Item(746, 5).Value = "White tail"
$ws.Cells.Item(746, 6).Value = "WT"
$ws.Cells.Item(747, 1).Value = 8
$ws.Cells.Item(747, 2).Formula = "=B738+1"
$ws.Cells.Item(747, 4).Value = 3
$ws.Cells.Item(747, 5).Value = "Half white tail"
$ws.Cells.Item(747, 6).Value = "1/2WT"
$ws.Cells.Item(748, 1).Value = 9
$ws.Cells.Item(748, 2).Formula = "=B747"
$ws.Cells.Item(748, 4).Value = 3
$ws.Cells.Item(748, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(748, 6).Value = "1/4WTT"
$ws.Cells.Item(749, 1).Value = 1
$ws.Cells.Item(749, 2).Formula = "=B740+1"
$ws.Cells.Item(749, 4).Value = 1
$ws.Cells.Item(749, 5).Value = "Black tipped tail"
$ws.Cells.Item(749, 6).Value = "BTT"
$ws.Cells.Item(750, 1).Value = 2
$ws.Cells.Item(750, 2).Formula = "=B749"
$ws.Cells.Item(750, 4).Value = 1
$ws.Cells.Item(750, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(750, 6).Value = "BDLRL"
$ws.Cells.Item(751, 1).Value = 3
$ws.Cells.Item(751, 2).Formula = "=B742+1"
$ws.Cells.Item(751, 4).Value = 1
$ws.Cells.Item(751, 5).Value = "White rear legs"
$ws.Cells.Item(751, 6).Value = "WRL"
$ws.Cells.Item(752, 1).Value = 4
$ws.Cells.Item(752, 2).Formula = "=B751"
$ws.Cells.Item(752, 4).Value = 2
$ws.Cells.Item(752, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(752, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(753, 1).Value = 5
$ws.Cells.Item(753, 2).Formula = "=B744+1"
$ws.Cells.Item(753, 4).Value = 2
$ws.Cells.Item(753, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(753, 6).Value = "1/2WT/WU"
$ws.Cells.Item(754, 1).Value = 6
$ws.Cells.Item(754, 2).Formula = "=B745+1"
$ws.Cells.Item(754, 4).Value = 2
$ws.Cells.Item(754, 5).Value = "Small white tip tail"
$ws.Cells.Item(754, 6).Value = "SWTT"
$ws.Cells.Item(755, 1).Value = 7
$ws.Cells.Item(755, 2).Formula = "=B754"
$ws.Cells.Item(755, 4).Value = 3
$ws.Cells.Item(755, 5).Value = "White tail"
$ws.Cells.Item(755, 6).Value = "WT"
$ws.Cells.Item(756, 1).Value = 8
$ws.Cells.Item(756, 2).Formula = "=B747+1"
$ws.Cells.Item(756, 4).Value = 3
$ws.Cells.Item(756, 5).Value = "Half white tail"
$ws.Cells.Item(756, 6).Value = "1/2WT"
$ws.Cells.Item(757, 1).Value = 9
$ws.Cells.Item(757, 2).Formula = "=B756"
$ws.Cells.Item(757, 4).Value = 3
$ws.Cells.Item(757, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(757, 6).Value = "1/4WTT"
$ws.Cells.Item(758, 1).Value = 1
$ws.Cells.Item(758, 2).Formula = "=B749+1"
$ws.Cells.Item(758, 4).Value = 1
$ws.Cells.Item(758, 5).Value = "Black tipped tail"
$ws.Cells.Item(758, 6).Value = "BTT"
$ws.Cells.Item(759, 1).Value = 2
$ws.Cells.Item(759, 2).Formula = "=B758"
$ws.Cells.Item(759, 4).Value = 1
$ws.Cells.Item(759, 5).Value = "Black dot left rear leg"
$ws.Cells.Item(759, 6).Value = "BDLRL"
$ws.Cells.Item(760, 1).Value = 3
$ws.Cells.Item(760, 2).Formula = "=B751+1"
$ws.Cells.Item(760, 4).Value = 1
$ws.Cells.Item(760, 5).Value = "White rear legs"
$ws.Cells.Item(760, 6).Value = "WRL"
$ws.Cells.Item(761, 1).Value = 4
$ws.Cells.Item(761, 2).Formula = "=B760"
$ws.Cells.Item(761, 4).Value = 2
$ws.Cells.Item(761, 5).Value = "Half white tail/Black spots under side"
$ws.Cells.Item(761, 6).Value = "1/2WT/BSU"
$ws.Cells.Item(762, 1).Value = 5
$ws.Cells.Item(762, 2).Formula = "=B753+1"
$ws.Cells.Item(762, 4).Value = 2
$ws.Cells.Item(762, 5).Value = "Half white tail/White under side"
$ws.Cells.Item(762, 6).Value = "1/2WT/WU"
$ws.Cells.Item(763, 1).Value = 6
$ws.Cells.Item(763, 2).Formula = "=B754+1"
$ws.Cells.Item(763, 4).Value = 2
$ws.Cells.Item(763, 5).Value = "Small white tip tail"
$ws.Cells.Item(763, 6).Value = "SWTT"
$ws.Cells.Item(764, 1).Value = 7
$ws.Cells.Item(764, 2).Formula = "=B763"
$ws.Cells.Item(764, 4).Value = 3
$ws.Cells.Item(764, 5).Value = "White tail"
$ws.Cells.Item(764, 6).Value = "WT"
$ws.Cells.Item(765, 1).Value = 8
$ws.Cells.Item(765, 2).Formula = "=B756+1"
$ws.Cells.Item(765, 4).Value = 3
$ws.Cells.Item(765, 5).Value = "Half white tail"
$ws.Cells.Item(765, 6).Value = "1/2WT"
$ws.Cells.Item(766, 1).Value = 9
$ws.Cells.Item(766, 2).Formula = "=B765"
$ws.Cells.Item(766, 4).Value = 3
$ws.Cells.Item(766, 5).Value = "One fourth white tip tail"
$ws.Cells.Item(766, 6).Value = "1/4WTT"

# --- Step 3: extend the AutoFilter range from A1:A676 to A1:A766 ---
# Toggle off then back on against the new range (AutoFilter() toggles state).
$ws.Range("A1:A766").AutoFilter() | Out-Null
$ws.Range("A1:A766").AutoFilter() | Out-Null

# --- Step 4: update the _FilterDatabase defined name to match the new range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Rat Colony!_FilterDatabase") {
        $n.RefersTo = "='Rat Colony'!`$A`$1:`$A`$766"
    }
}

# --- Step 5: move the active selection to match the edited workbook (B670) ---
$ws.Range("B670").Select()
